# Adding Self service user testcases into Logixal QA box 2 env
# Appends two new rows (32, 33) to the MasterExecutor sheet for:
#   TC42_Verify_PlaceOrder_SelfService_SingleUser
#   TC43_Verify_PlaceOrder_SelfService_MultiUser

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clone the formatting (borders / wrap / alignment / font) of an existing
# "ht=30" data row (row 3) onto the two new rows, without disturbing any of
# the workbook's existing styles.
$ws.Range("A3:F3").Copy()
$ws.Range("A32:F32").PasteSpecial(-4122)
$ws.Range("A3:F3").Copy()
$ws.Range("A33:F33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 32: TC42_Verify_PlaceOrder_SelfService_SingleUser
$ws.Cells.Item(32, 1).Value2 = "ALL_PAGES"
$ws.Cells.Item(32, 2).Value2 = "END_TO_END"
$ws.Cells.Item(32, 3).Value2 = "TC42_Verify_PlaceOrder_SelfService_SingleUser"

# Row 33: TC43_Verify_PlaceOrder_SelfService_MultiUser
$ws.Cells.Item(33, 1).Value2 = "ALL_PAGES"
$ws.Cells.Item(33, 2).Value2 = "END_TO_END"
$ws.Cells.Item(33, 3).Value2 = "TC43_Verify_PlaceOrder_SelfService_MultiUser"

# Descriptions (column D) filled in after both testcase numbers, matching
# the order the new shared strings were authored in.
$ws.Cells.Item(32, 4).Value2 = "Verify place order using Selfservice Single user"
$ws.Cells.Item(33, 4).Value2 = "Verify place order using Selfservice multi user"

# RunMode / Severity columns for both new rows.
$ws.Cells.Item(32, 5).Value2 = "Yes"
$ws.Cells.Item(32, 6).Value2 = "High"
$ws.Cells.Item(33, 5).Value2 = "Yes"
$ws.Cells.Item(33, 6).Value2 = "High"

# Match the taller row height used by other multi-requirement rows.
$ws.Rows.Item(32).RowHeight = 30
$ws.Rows.Item(33).RowHeight = 30

# Update the visible selection to reflect the newly extended column E range.
$ws.Range("E2:E33").Select()
